# "actualizacion de combobox terminado"
#
# Updates the "Sheet" control table (Input Mode / No.Variables / No.ItemsCoeffM
# / ItemsCoeffM) to 23 variables with 3 coefficient items, appends the
# "Coord *" header/value rows describing the index ranges used by each of the
# coefficient-matrix worksheets, grows the "diffusion" matrix from 22x22 to
# 23x23 (updating its first-row coefficients), and materialises the
# "absorption" (23x23 matrix) and "source" (23x1 vector) worksheets that were
# still empty placeholders.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sheet": control/summary table
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Sheet")

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 23
$ws.Range("C2").Value = 3
$ws.Range("D2").Value = "1,2,3"

$ws.Range("A3").Value = "Coord Diffusion"
$ws.Range("B3").Value = "Coord Absorption"
$ws.Range("C3").Value = "Coord Source"
$ws.Range("D3").Value = "Coord Mass"
$ws.Range("E3").Value = "Coord DamMass"
$ws.Range("F3").Value = "Coord CFlux"
$ws.Range("G3").Value = "Coord Convection"
$ws.Range("H3").Value = "Coord CSource"

$ws.Range("A4").Value = "[0, 1]"
$ws.Range("B4").Value = "[0, 1]"
$ws.Range("C4").Value = "[0]"
$ws.Range("D4").Value = "[0, 0]"
$ws.Range("E4").Value = "[0, 0]"
$ws.Range("F4").Value = "[0, 0]"
$ws.Range("G4").Value = "[0, 0]"
$ws.Range("H4").Value = "[0]"

# ---------------------------------------------------------------------
# Sheet "diffusion": grow the 22x22 matrix to 23x23 and tweak row 1
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("diffusion")

# Touch the full new extent (with a text format) so every cell in the
# bigger 23x23 block is materialised - matching the exported worksheet's
# fully-populated (if blank) grid - then drop the formatting again so no
# stray cell styles are left behind.
$ws.Range("A1:W23").NumberFormat = "@"
$ws.Range("B1").Value = "[45, 0, 0, 45]"
$ws.Range("C1").ClearContents()
$ws.Cells.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "absorption": new 23x23 matrix (first row holds two settings)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("absorption")

$ws.Range("A1:W23").NumberFormat = "@"
$ws.Range("A1").Value = "22"
$ws.Range("B1").Value = "10"
$ws.Cells.ClearFormats()

# ---------------------------------------------------------------------
# Sheet "source": new 23x1 vector (first cell holds a setting)
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("source")

$ws.Range("A1:A23").NumberFormat = "@"
$ws.Range("A1").Value = "5"
$ws.Cells.ClearFormats()
